# Generate Report for Handback
# Update the handback status + datetime for c33b6286-7556-4d32-9fc9-7304f58620f1.md
# across the Overview sheet and each locale (zh-cn, de-de) detail sheet, clearing
# the stale "handback file not latest" error now that the handback succeeded.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the c33b6286... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn detail sheet: row 3 is the c33b6286... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-21 22:56:59"
$wsZhCn.Range("P3").Value = ""

# --- de-de detail sheet: row 3 is the c33b6286... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-21 22:57:10"
$wsDeDe.Range("P3").Value = ""
